$wb = $excel.ActiveWorkbook

# Station 2
$ws = $wb.Worksheets.Item("Station 2")
$ws.Range("C2").Value = 2.579353533977641
$ws.Range("D2").Value = 3.1166482402283155
$ws.Range("C3").Value = 1.2703958726453837
$ws.Range("D3").Value = 0.22734744821688707
$ws.Range("C4").Value = 0.9361510744866999
$ws.Range("D4").Value = 0.2685088269917376
$ws.Range("C5").Value = 0.7646300521535654
$ws.Range("D5").Value = 0.2743406468987301
$ws.Range("C6").Value = 0.6537399203573111
$ws.Range("D6").Value = 0.26838821583734984
$ws.Range("C7").Value = 0.572875358122322
$ws.Range("D7").Value = 0.26157347014559945
$ws.Range("C8").Value = 0.509415203741657
$ws.Range("D8").Value = 0.26837254599920757
$ws.Range("C9").Value = 0.4571404618118016
$ws.Range("D9").Value = 0.26915309525570946
$ws.Range("C10").Value = 0.41257238112573447
$ws.Range("D10").Value = 0.23046800642364804
$ws.Range("C11").Value = 0.37355608568020626
$ws.Range("D11").Value = 0.21116645656476302
$ws.Range("C12").Value = 0.3386406516432431
$ws.Range("D12").Value = 0.16118981500489993
$ws.Range("C13").Value = 0.30677869695506654
$ws.Range("D13").Value = 0.0778053663493638
$ws.Range("C14").Value = 0.27716423268596124
$ws.Range("D14").Value = 0.025681244695200812
$ws.Range("C15").Value = 0.24913305677524808
$ws.Range("D15").Value = 0.06150676461412733
$ws.Range("C16").Value = 0.22208976721304097
$ws.Range("D16").Value = 0.07162950510489025
$ws.Range("C17").Value = 0.19543843735453054
$ws.Range("D17").Value = 0.1865182495580229
$ws.Range("C18").Value = 0.16849018645113747
$ws.Range("D18").Value = 0.07938023441502513
$ws.Range("C19").Value = 0.14028834237199356
$ws.Range("D19").Value = 0.001456604170433733
$ws.Range("C20").Value = 0.10914767353926497
$ws.Range("D20").Value = 0.14558438209902969
$ws.Range("C21").Value = 0.07082181360114208
$ws.Range("D21").Value = 8.403144603104638

# Station 4
$ws = $wb.Worksheets.Item("Station 4")
$ws.Range("C2").Value = 2.9601802826755814
$ws.Range("D2").Value = 1.1240854225797086
$ws.Range("C3").Value = 1.4396033021221974
$ws.Range("D3").Value = 0.2636508018471132
$ws.Range("C4").Value = 1.047977664319615
$ws.Range("D4").Value = 0.32070175030084397
$ws.Range("C5").Value = 0.8465821949580427
$ws.Range("D5").Value = 0.33344988574326645
$ws.Range("C6").Value = 0.7171284583406702
$ws.Range("D6").Value = 0.36011474989835257
$ws.Range("C7").Value = 0.6238528920754581
$ws.Range("D7").Value = 0.37235232311740735
$ws.Range("C8").Value = 0.5517259893358106
$ws.Range("D8").Value = 0.37310969643262004
$ws.Range("C9").Value = 0.4931596553159977
$ws.Range("D9").Value = 0.3718465011372181
$ws.Range("C10").Value = 0.44383477734482
$ws.Range("D10").Value = 0.37559341387555406
$ws.Range("C11").Value = 0.40106946786266523
$ws.Range("D11").Value = 0.3557486671084613
$ws.Range("C12").Value = 0.363078947790894
$ws.Range("D12").Value = 0.35205788973597835
$ws.Range("C13").Value = 0.32860043394652594
$ws.Range("D13").Value = 0.3396940580445786
$ws.Range("C14").Value = 0.29668411959775953
$ws.Range("D14").Value = 0.3067130275748174
$ws.Range("C15").Value = 0.2665643745249468
$ws.Range("D15").Value = 0.251942066690424
$ws.Range("C16").Value = 0.23756832718548854
$ws.Range("D16").Value = 0.18422438407772485
$ws.Range("C17").Value = 0.2090340977252289
$ws.Range("D17").Value = 0.12928518707856262
$ws.Range("C18").Value = 0.18020811699424072
$ws.Range("D18").Value = 0.0036293211388430816
$ws.Range("C19").Value = 0.15005687595939374
$ws.Range("D19").Value = 0.24454114866739868
$ws.Range("C20").Value = 0.11677348025632082
$ws.Range("D20").Value = 1.0275836726605585
$ws.Range("C21").Value = 0.07581050306706753
$ws.Range("D21").Value = 9.065502230547906

# Station 6
$ws = $wb.Worksheets.Item("Station 6")
$ws.Range("C2").Value = 3.124402630899344
$ws.Range("D2").Value = 0.6123167336492975
$ws.Range("C3").Value = 1.5220509481606097
$ws.Range("D3").Value = 0.2638900126626336
$ws.Range("C4").Value = 1.1092300053091255
$ws.Range("D4").Value = 0.2783511561623633
$ws.Range("C5").Value = 0.8965257785442293
$ws.Range("D5").Value = 0.30675622938149705
$ws.Range("C6").Value = 0.7593801566775402
$ws.Range("D6").Value = 0.3383024280641013
$ws.Range("C7").Value = 0.6601900407749585
$ws.Range("D7").Value = 0.34805549153248155
$ws.Range("C8").Value = 0.5831910839486772
$ws.Range("D8").Value = 0.3422817943595642
$ws.Range("C9").Value = 0.520444058137637
$ws.Range("D9").Value = 0.3326349248040764
$ws.Range("C10").Value = 0.46744088632773634
$ws.Range("D10").Value = 0.33341815281434095
$ws.Range("C11").Value = 0.4213893377338488
$ws.Range("D11").Value = 0.327601255893088
$ws.Range("C12").Value = 0.3804366952420054
$ws.Range("D12").Value = 0.3417951982140664
$ws.Range("C13").Value = 0.3432772813266313
$ws.Range("D13").Value = 0.3569750202461287
$ws.Range("C14").Value = 0.308934439118261
$ws.Range("D14").Value = 0.3273689470525396
$ws.Range("C15").Value = 0.27662668847150523
$ws.Range("D15").Value = 0.273904189629524
$ws.Range("C16").Value = 0.2456730446249153
$ws.Range("D16").Value = 0.15792553282397392
$ws.Range("C17").Value = 0.21540830916596854
$ws.Range("D17").Value = 0.007033374267451251
$ws.Range("C18").Value = 0.1850762152388011
$ws.Range("D18").Value = 0.16452686691388707
$ws.Range("C19").Value = 0.15363345329886327
$ws.Range("D19").Value = 0.3788001259190506
$ws.Range("C20").Value = 0.1192405244148076
$ws.Range("D20").Value = 1.237559280691539
$ws.Range("C21").Value = 0.07725273397025742
$ws.Range("D21").Value = 9.256989924016718

# Station 8
$ws = $wb.Worksheets.Item("Station 8")
$ws.Range("C2").Value = 2.823066759602783
$ws.Range("D2").Value = 0.6414821077469524
$ws.Range("C3").Value = 1.3427766958024134
$ws.Range("D3").Value = 0.2384251576046014
$ws.Range("C4").Value = 0.9507786165541559
$ws.Range("D4").Value = 0.29544917964166467
$ws.Range("C5").Value = 0.7443012268939574
$ws.Range("D5").Value = 0.27653191726058063
$ws.Range("C6").Value = 0.609377199628987
$ws.Range("D6").Value = 0.3149352032636976
$ws.Range("C7").Value = 0.5113840049947898
$ws.Range("D7").Value = 0.31995275064206474
$ws.Range("C8").Value = 0.43574580906165084
$ws.Range("D8").Value = 0.3003372864207985
$ws.Range("C9").Value = 0.37506932473430293
$ws.Range("D9").Value = 0.3046079963980011
$ws.Range("C10").Value = 0.32509018977681625
$ws.Range("D10").Value = 0.2788878029039363
$ws.Range("C11").Value = 0.2830986932366888
$ws.Range("D11").Value = 0.2659765040303097
$ws.Range("C12").Value = 0.24723870789777377
$ws.Range("D12").Value = 0.2579963002661648
$ws.Range("C13").Value = 0.21616191615593516
$ws.Range("D13").Value = 0.2387174447442414
$ws.Range("C14").Value = 0.18884119942829702
$ws.Range("D14").Value = 0.24391238385419947
$ws.Range("C15").Value = 0.1644597085032884
$ws.Range("D15").Value = 0.18526923543004525
$ws.Range("C16").Value = 0.1423360795552088
$ws.Range("D16").Value = 0.1684779944112836
$ws.Range("C17").Value = 0.12186336076279003
$ws.Range("D17").Value = 0.17911941290056568
$ws.Range("C18").Value = 0.10244179024130008
$ws.Range("D18").Value = 0.6693914137722262
$ws.Range("C19").Value = 0.08336958828623307
$ws.Range("D19").Value = 0.9361855238176288
$ws.Range("C20").Value = 0.06357565615308375
$ws.Range("D20").Value = 2.579008901953164
$ws.Range("C21").Value = 0.04057725359862088
$ws.Range("D21").Value = 4.387517825136961
